$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header cell C1 "Gender" -> "Status"
$ws.Range("C1").Value = "Status"

# Update the notes cell E3 text
$ws.Range("E3").Value = "*DO NOT DELETE TITLE ROW`n- Status is either 1 or 0 (1 is true, 0 is false)`n- Email is either @gmail.com or @fpt.edu.vn`n- DO NOT make gmail into a url/link, it will bug"

# Change selection from F3 to E3
$ws.Range("E3").Select()

# Change row 3 height from 90 to 105
$ws.Rows("3:3").RowHeight = 105
